$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header cells ---
$ws.Range("A1").Value = "Url"
$ws.Range("B1").Value = "query"
# C1 "dbExcel" and D1 "WebExcel" are unchanged.

# --- Row 2 data cells ---
# A2 becomes the Bento QA URL, styled + linked as a hyperlink (style reset to
# Normal first so the Hyperlink style doesn't inherit the old wrapText align).
$ws.Range("A2").Style = "Normal"
$ws.Range("A2").Value = "https://bento-qa.bento-tools.org/"
$ws.Hyperlinks.Add($ws.Range("A2"), "https://bento-qa.bento-tools.org/")

# B2 keeps its word-wrap style but gets the updated Cypher query text
# (er_status flipped to "Positive"; endocrine_therapy_type / head(labels(samp))
# clauses dropped).
$ws.Range("B2").Value = "MATCH (ss:study_subject)
MATCH (samp)-[:sample_of_study_subject]->(ss)
MATCH (ss)<-[:sample_of_study_subject]-(samp:sample)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)
WITH DISTINCT ss, samp, collect(DISTINCT samp.sample_id) AS samples, collect(DISTINCT lp.laboratory_procedure_id) AS lab_procedures, collect(DISTINCT f) AS files
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
MATCH (ss)<-[:demographic_of_study_subject]-(demo)
WHERE ss.disease_subtype IN [`"Tubular Carcinoma`"] and d.tumor_size_group In [`"(3,4]`"] and d.er_status In [`"Positive`"]and d.pr_status In [`"Positive`"] 
return DISTINCT ss.study_subject_id as ``Case ID``,
   p.program_acronym as ``Program Code``,
    p.program_id as Program_ID,
   s.study_acronym as ``Arm``,
   ss.disease_subtype as ``Diagnosis``,
   sf.grouped_recurrence_score AS ``Recurrence Score``,
   d.tumor_size_group AS ``tumor_size``,
   d.er_status AS ``ER Status``,
   d.pr_status AS ``PR Status``,
   demo.age_at_index AS ``Age (years)``,
	demo.survival_time AS ``Survival (days)``"

# C2 / D2: swap in the manifest filename, web-data filename stays the same text.
$ws.Range("C2").Value = "TC03_Bento_E2E_Select-Single-CaseDetail_Manifest.xlsx"
$ws.Range("D2").Value = "TC03_Bento_E2E_Select-Single-CaseDetail_WebData.xlsx"

# Row 2 height shrinks now that the query text is a bit shorter.
$ws.Rows.Item(2).RowHeight = 375

# Selection moves off the corner cell.
$ws.Range("D7").Select() | Out-Null
